$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date for all existing data rows 2-447
for ($r = 2; $r -le 447; $r++) {
    $ws.Cells.Item($r, 3).Value = 45203
}

# Row 447 gains an explicit row height (matches default visual row height)
$ws.Rows.Item(447).RowHeight = 15

# Add new row 448
$ws.Range("B448:C448").NumberFormat = "YYYY-MM-DD"
$ws.Range("A448").Value = "A 47034-2023"
$ws.Range("B448").Value = 45201
$ws.Range("C448").Value = 45203
$ws.Range("D448").Value = "UPPSALA LÄN"
$ws.Range("E448").Value = "TIERP"
$ws.Range("F448").Value = "Bergvik skog öst AB"
$ws.Range("G448").Value = 4.2
$ws.Range("H448").Value = 0
$ws.Range("I448").Value = 0
$ws.Range("J448").Value = 0
$ws.Range("K448").Value = 0
$ws.Range("L448").Value = 0
$ws.Range("M448").Value = 0
$ws.Range("N448").Value = 0
$ws.Range("O448").Value = 0
$ws.Range("P448").Value = 0
$ws.Range("Q448").Value = 0
$ws.Range("R448").WrapText = $true
